$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.414.52'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.863.23'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -1.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.69'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E6").Value = '  -1.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5066'
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08313'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.38'
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.100'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").Value = '1.867.77'
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.24'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.218'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("E16").Value = '  -1.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.99'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06717'
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.56'
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.892'
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").Value = '28.462.18'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.198'
$ws.Range("E25").Value = '  -4.34%  '
$ws.Range("D26").Value = '2.079.93'
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.90'
$ws.Range("E27").Value = '  -3.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.58'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.47'
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1035'
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.738'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.605'
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02438'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06550'
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.921'
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2149'
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.001'
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.175'
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.231'
$ws.Range("E41").Value = '  -3.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6326'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.05'
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.007'
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5968'
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.99'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.678'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.988'
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.06'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.203'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.146'
$ws.Range("E51").Value = '  -6.69%  '
